$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating the existing "2022-Q2"
#    sheet (so sheetPr/styles/page setup all come along for free), placed
#    right before it. Tab order ends up: 总计, 2022-Q3, 2022-Q2, 2022-Q1.
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item(2)
$wsQ2.Copy($wsQ2)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Trim the duplicated sheet down to just the 2 data rows the Q3 table needs
# (it currently has the 15 rows copied from 2022-Q2, rows 4-16).
for ($r = 16; $r -ge 4; $r--) {
    $wsQ3.Rows.Item($r).Delete()
}

# Overwrite the two data rows with the actual 2022-Q3 fund holdings.
# Columns B-G are text cells (even the numeric-looking ones), so force a
# text number format before writing, then clear it again so the cells stay
# unstyled like the rest of the table.
$wsQ3.Range("B2:G3").NumberFormat = "@"

$wsQ3.Range("B2").Value = "013920"
$wsQ3.Range("C2").Value = "兴华创新医疗6个月持有混合A"
$wsQ3.Range("D2").Value = "0.18"
$wsQ3.Range("E2").Value = "94.83"
$wsQ3.Range("F2").Value = "6.47"
$wsQ3.Range("G2").Value = "0.0116"
$wsQ3.Range("H2").Value = 5

$wsQ3.Range("B3").Value = "013921"
$wsQ3.Range("C3").Value = "兴华创新医疗6个月持有混合C"
$wsQ3.Range("D3").Value = "0.05"
$wsQ3.Range("E3").Value = "94.83"
$wsQ3.Range("F3").Value = "6.47"
$wsQ3.Range("G3").Value = "0.0032"
$wsQ3.Range("H3").Value = 5

$wsQ3.Range("B2:G3").ClearFormats()

# ---------------------------------------------------------------------------
# 2. Insert a new summary row for "2022-Q3" at the top of the "总计" sheet's
#    data (row 2), pushing the existing 2022-Q2 / 2022-Q1 rows down.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.01

# Renumber the index column for the rows that shifted down.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
